$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cellRef, $value)
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $ws.Range("C2").Copy()
    $cell.PasteSpecial(-4122)
}

Set-TextCell "E2" "2026-02-23 17:48:42"
Set-TextCell "H2" "34%"
Set-TextCell "E3" "2026-02-23 17:48:45"
Set-TextCell "H3" "32%"
Set-TextCell "K3" "16.7 MJ/m2"
Set-TextCell "E4" "2026-02-23 17:48:48"
Set-TextCell "J4" "1024.8 hPa"
Set-TextCell "K4" "14.9 MJ/m2"
Set-TextCell "O4" "12.5 °C"
Set-TextCell "E5" "2026-02-23 17:48:50"
Set-TextCell "K5" "16.4 MJ/m2"
Set-TextCell "E6" "2026-02-23 17:48:53"
Set-TextCell "O6" "14.1 °C"
Set-TextCell "E7" "2026-02-23 17:48:56"
Set-TextCell "O7" "14.2 °C"
Set-TextCell "E8" "2026-02-23 17:48:58"
Set-TextCell "J8" "1024.4 hPa"
Set-TextCell "K8" "15.9 MJ/m2"
Set-TextCell "E9" "2026-02-23 17:49:01"
Set-TextCell "H9" "70%"
Set-TextCell "O9" "12.8 °C"
Set-TextCell "E10" "2026-02-23 17:49:04"
Set-TextCell "O10" "11.3 °C"
Set-TextCell "E11" "2026-02-23 17:49:06"
Set-TextCell "O11" "9.0 °C"
Set-TextCell "E12" "2026-02-23 17:49:09"
Set-TextCell "O12" "10.9 °C"
Set-TextCell "E13" "2026-02-23 17:49:12"
Set-TextCell "H13" "59%"
Set-TextCell "J13" "1027.1 hPa"
Set-TextCell "L13" "24.8 km/h - 110º 17:18 TU"
Set-TextCell "O13" "6.9 °C"
Set-TextCell "E14" "2026-02-23 17:49:14"
Set-TextCell "E15" "2026-02-23 17:49:17"
Set-TextCell "O15" "13.0 °C"
Set-TextCell "E16" "2026-02-23 17:49:19"
Set-TextCell "E17" "2026-02-23 17:49:22"
Set-TextCell "K17" "17.3 MJ/m2"
Set-TextCell "E18" "2026-02-23 17:49:25"
Set-TextCell "H18" "72%"
Set-TextCell "O18" "11.0 °C"
Set-TextCell "E19" "2026-02-23 17:49:28"
Set-TextCell "K19" "15.4 MJ/m2"
Set-TextCell "E20" "2026-02-23 17:49:30"
Set-TextCell "K20" "16.7 MJ/m2"
Set-TextCell "E21" "2026-02-23 17:49:33"
Set-TextCell "H21" "59%"
Set-TextCell "J21" "1026.0 hPa"
Set-TextCell "O21" "9.5 °C"
Set-TextCell "E22" "2026-02-23 17:49:36"
Set-TextCell "E23" "2026-02-23 17:49:38"
Set-TextCell "H23" "20%"
Set-TextCell "K23" "16.4 MJ/m2"
Set-TextCell "O23" "3.7 °C"
Set-TextCell "E24" "2026-02-23 17:49:41"
Set-TextCell "H24" "82%"
Set-TextCell "J24" "1026.4 hPa"
Set-TextCell "K24" "16.2 MJ/m2"
Set-TextCell "O24" "8.5 °C"
Set-TextCell "E25" "2026-02-23 17:49:44"
Set-TextCell "K25" "17.2 MJ/m2"
Set-TextCell "O25" "6.1 °C"
Set-TextCell "E26" "2026-02-23 17:49:46"
Set-TextCell "J26" "1023.7 hPa"
Set-TextCell "E27" "2026-02-23 17:49:49"
Set-TextCell "H27" "27%"
Set-TextCell "E28" "2026-02-23 17:49:52"
Set-TextCell "H28" "66%"
Set-TextCell "O28" "11.2 °C"
Set-TextCell "E29" "2026-02-23 17:49:54"
Set-TextCell "H29" "80%"
Set-TextCell "K29" "15.7 MJ/m2"
Set-TextCell "O29" "10.9 °C"
Set-TextCell "E30" "2026-02-23 17:49:57"
Set-TextCell "H30" "67%"
Set-TextCell "E31" "2026-02-23 17:50:00"
Set-TextCell "O31" "16.4 °C"
Set-TextCell "E32" "2026-02-23 17:50:02"
Set-TextCell "H32" "64%"
Set-TextCell "K32" "16.0 MJ/m2"
Set-TextCell "O32" "8.5 °C"
Set-TextCell "E33" "2026-02-23 17:50:05"
Set-TextCell "J33" "1025.5 hPa"
Set-TextCell "O33" "8.6 °C"
Set-TextCell "E34" "2026-02-23 17:50:08"
Set-TextCell "H34" "38%"
Set-TextCell "O34" "4.5 °C"
Set-TextCell "E35" "2026-02-23 17:50:10"
Set-TextCell "K35" "16.9 MJ/m2"
Set-TextCell "E36" "2026-02-23 17:50:13"
Set-TextCell "K36" "15.3 MJ/m2"
Set-TextCell "O36" "13.0 °C"
Set-TextCell "E37" "2026-02-23 17:50:16"
Set-TextCell "J37" "1026.6 hPa"
Set-TextCell "O37" "9.6 °C"
Set-TextCell "E38" "2026-02-23 17:50:18"
Set-TextCell "O38" "12.4 °C"
Set-TextCell "E39" "2026-02-23 17:50:21"
Set-TextCell "H39" "22%"
Set-TextCell "E40" "2026-02-23 17:50:23"
Set-TextCell "J40" "1026.3 hPa"
Set-TextCell "O40" "9.0 °C"
Set-TextCell "E41" "2026-02-23 17:50:26"
Set-TextCell "J41" "1024.8 hPa"
Set-TextCell "K41" "15.9 MJ/m2"
Set-TextCell "E42" "2026-02-23 17:50:29"
Set-TextCell "H42" "76%"
Set-TextCell "O42" "11.9 °C"
Set-TextCell "E43" "2026-02-23 17:50:31"
Set-TextCell "H43" "70%"
Set-TextCell "O43" "10.2 °C"
Set-TextCell "E44" "2026-02-23 17:50:34"
Set-TextCell "K44" "16.3 MJ/m2"
Set-TextCell "E45" "2026-02-23 17:50:37"
Set-TextCell "K45" "14.5 MJ/m2"
Set-TextCell "E46" "2026-02-23 17:50:40"
Set-TextCell "H46" "73%"
Set-TextCell "J46" "1026.3 hPa"
Set-TextCell "O46" "10.1 °C"
